$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for rows 2-7 from 45224 (2023-10-25)
# to 45233 (2023-11-03), preserving existing cell formatting.
$ws.Range("C2:C7").Value = 45233
